# Fruta / hortaliza, semanal
# Two new weekly price rows are inserted at the top of the data block
# (right after the header row), pushing the existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 7 (the first data row for
# this market/product) so the historical rows shift down to 9..14.
$ws.Rows("7:8").Insert()

# Columns A,B,C,E,F,G,H,I,J,K are constant for every record in this sheet
# (same market / product taxonomy) - copy them down from the row that is
# now at position 9 (the former row 7) into the two freshly inserted rows.
foreach ($col in @("A","B","C","E","F","G","H","I","J","K")) {
    $value = $ws.Range("$col`9").Value2
    $ws.Range("$col`7").Value = $value
    $ws.Range("$col`8").Value = $value
}

# Row 7 - new "Primera" quality entry
$ws.Range("D7").Value = 44904
$ws.Range("D7").NumberFormat = $ws.Range("D9").NumberFormat
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 16000
$ws.Range("P7").Value = 15500
$ws.Range("Q7").Value = "`$/bandeja 10 kilos"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1550
$ws.Range("T7").Value = 10

# Row 8 - new "Segunda" quality entry
$ws.Range("D8").Value = 44904
$ws.Range("D8").NumberFormat = $ws.Range("D9").NumberFormat
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 14000
$ws.Range("Q8").Value = "`$/bandeja 10 kilos"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 1400
$ws.Range("T8").Value = 10
